$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.443.97'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '2.021.58'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'254.26"
$ws.Range("E5").Value = '  +3.88%  '
$ws.Range("D6").Value = "'0.619"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'56.85"
$ws.Range("E8").Value = '  -6.81%  '
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = "'0.102"
$ws.Range("D12").Value = "'14.51"
$ws.Range("E12").Value = '  -2.73%  '
$ws.Range("D13").Value = '2.324.37'
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = "'0.816"
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = "'21.12"
$ws.Range("E15").Value = '  -5.13%  '
$ws.Range("D16").Value = "'5.34"
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").Value = '2.032.74'
$ws.Range("E17").Value = '  +2.58%  '
$ws.Range("D18").Value = '37.375.65'
$ws.Range("E18").Value = '  +1.38%  '
$ws.Range("D19").Value = "'69.45"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").Value = '0.0₃0848'
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").Value = "'5.17"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = "'228.42"
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = "'2.60"
$ws.Range("E24").Value = '  +3.23%  '
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").Value = "'164.03"
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").Value = "'9.03"
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").Value = "'19.86"
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("E29").Value = '  -11.64%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = "'0.0662"
$ws.Range("E32").Value = '  +6.67%  '
$ws.Range("D33").Value = "'4.70"
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = '  +5.41%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("E39").Value = '  -3.07%  '
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").Value = "'0.0964"
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("E43").Value = '  +1.11%  '
$ws.Range("D44").Value = '1.394.70'
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("D45").Value = "'15.99"
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").Value = "'90.67"
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = "'7.34"
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("D51").Value = '2.215.77'
$ws.Range("E51").Value = '  +2.46%  '
